$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.623558044433594
$ws.Range("B1").Value = 1.501290678977966
$ws.Range("C1").Value = 2.000089645385742
$ws.Range("D1").Value = 1.796268939971924
$ws.Range("E1").Value = 2.869723796844482
